$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.989.03'
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').Value = '1.796.32'
$ws.Range('E3').Value = '  -2.70%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.28'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  -2.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3601'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07259'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8409'
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.26'
$ws.Range('E11').Value = '  -3.90%  '
$ws.Range('D12').Value = '1.823.89'
$ws.Range('E12').Value = '  -3.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.282'
$ws.Range('E13').Value = '  -3.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.367'
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06778'
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '80.54'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008742'
$ws.Range('E18').Value = '  -3.62%  '
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.98'
$ws.Range('E20').Value = '  -4.03%  '
$ws.Range('D21').Value = '27.147.55'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.083'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.04'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = '2.054.43'
$ws.Range('E24').Value = '  -3.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.954'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.98'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.13'
$ws.Range('E27').Value = '  -4.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.004'
$ws.Range('E28').Value = '  -6.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.47'
$ws.Range('E29').Value = '  -1.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.650'
$ws.Range('E30').Value = '  -12.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08993'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7277'
$ws.Range('E32').Value = '  -8.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.861'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.335'
$ws.Range('E34').Value = '  -6.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.091'
$ws.Range('E35').Value = '  -7.04%  '
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.080'
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05146'
$ws.Range('E38').Value = '  -5.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01902'
$ws.Range('E39').Value = '  -3.17%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1630'
$ws.Range('E40').Value = '  -3.97%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4974'
$ws.Range('E41').Value = '  -4.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.619'
$ws.Range('E42').Value = '  -7.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.068'
$ws.Range('E43').Value = '  -6.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.947'
$ws.Range('E44').Value = '  -12.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.02'
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.20'
$ws.Range('E46').Value = '  -4.39%  '
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06307'
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4520'
$ws.Range('E49').Value = '  -5.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.598'
$ws.Range('E50').Value = '  -4.07%  '
$ws.Range('E51').Value = '  -7.65%  '
